# Fruta / hortaliza, semanal
# Insert a new weekly record as row 102, pushing all subsequent records down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102 (shifts existing rows 102..207 down to 103..208)
$ws.Rows("102:102").Insert()

# Populate the newly inserted row with the new data point
$ws.Cells.Item(102, 1).Value  = 4
$ws.Cells.Item(102, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(102, 3).Value  = 'Los Lagos'
$ws.Cells.Item(102, 4).Value  = 44705
$ws.Cells.Item(102, 5).Value  = 10
$ws.Cells.Item(102, 6).Value  = 'Fruta'
$ws.Cells.Item(102, 7).Value  = 100108
$ws.Cells.Item(102, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(102, 9).Value  = 100108002
$ws.Cells.Item(102, 10).Value = 'Mango'
$ws.Cells.Item(102, 11).Value = 'Sin especificar'
$ws.Cells.Item(102, 12).Value = 'Primera'
$ws.Cells.Item(102, 13).Value = 200
$ws.Cells.Item(102, 14).Value = 7500
$ws.Cells.Item(102, 15).Value = 8000
$ws.Cells.Item(102, 16).Value = 7750
$ws.Cells.Item(102, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(102, 18).Value = 'Perú'
$ws.Cells.Item(102, 19).Value = 1938
$ws.Cells.Item(102, 20).Value = 4
